$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 64-73 currently only have a (blank) styled cell in column A.
# Extend the existing data pattern seen in rows 4-63:
#   column A = sequential depth index (continues 58,59,60,61 -> 61..70)
#   column C (AGE) = 0
#   column D (zbio) = 1
#   column E (ABU)  = 0
#   column F (ISO)  = 3
for ($i = 0; $i -le 9; $i++) {
    $row = 64 + $i
    $ws.Cells.Item($row, 1).Value = 61 + $i   # column A
    $ws.Cells.Item($row, 3).Value = 0          # column C
    $ws.Cells.Item($row, 4).Value = 1          # column D
    $ws.Cells.Item($row, 5).Value = 0          # column E
    $ws.Cells.Item($row, 6).Value = 3          # column F
}

# Match the updated view/selection state: scrolled so row 31 is visible
# at the top, with A64:F73 selected (active cell A64).
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A64:F73").Select()
